$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G-column comment reformatting: add "Comentarios personalizados:" header
# and convert each "> " bullet line into a "- " bullet line (44 rows). ---
$ws.Range("G17").Value = "> Comentarios personalizados:`n- La idea que estás desarrollando responde a la cotidianidad de sus habitantes y se relaciona de manera adecuada con el entorno. `n- Debes hacer un esfuerzo por elevar el nivel de representación, pues nos costó mucho diferenciar los diferentes elementos que componen el espacio. `n- Debido a lo anterior, no nos fue posible verificar qué tan bien resueltas están las zonas de circulación, la proporción de los espacios con relación al uso que tienen, la manera como iluminas y ventilas los espacios de alta permanencia y la manera como resuelves la cubierta.`n- No presentaste cuadro de áreas, ni planimetría en escala 1:250`n- La nota global de tu entrega es 4,1, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 3,1 (rebaja de 1 unidad)."
$ws.Range("G42").Value = "> Comentarios personalizados:`n- El proyecto interpreta la cotidianidad de sus habitantes y se relaciona de manera adecuada con su entorno. Asi mismo, se resuelve de manera adecuada la ventilación e iluminación de los diferentes espacios y el plantemiento de la cubierta.`n- Los espacios de circulación y las proporciones de los espacios según su uso, están bien resueltos de manera general; es posible no obstante, refinar un poco más esta parte. `n- Aunque presentaste a tiempo todos los productos requeridos para la entrega a través del set de planos, la memoria solo apareció en un momento extemporáneo.`n- La nota global de tu entrega es 4,3."
$ws.Range("G67").Value = "> Comentarios personalizados:`n-  El proyecto interpreta la cotidianidad de sus habitantes y se relaciona de manera adecuada con su entorno. `n- Muy interesante que te lances a proponer cubiertas abovedadas: esto genera una espacialidad singular con relación al resto de proyectos del taller y puede funcionar muy bien para ese clima.`n- De manera general, estás empezando a resolver de manera adecuada la ventilación e iluminación de los diferentes espacios, así como la proporción que tienen con relación a su uso y el tamaño de las áreas de circulación; no obstante, es posibe hacer ajustes en este aspecto.`n- La nota global de tu entrega es 4,3, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 3,3 (rebaja de 1 unidad)."
$ws.Range("G92").Value = "> Comentarios personalizados:`n- Los espacios que estás considerando tienen buenas posibiidades de iluminarse y ventilarse.`n- La manera como resuelves los espacios de circulación, así como el tamaño de los espacios con relación a su uso, admiten algunos ajustes para resolverlos un poco mejor. `n- En las cubiertas es importante explorar alternativas pues no vemos una solución adecuada en este momento. También, revisa la manera como estás interpretando la cotidianidad de los habitantes de esta vivienda.`n- Te retamos a elevar tu nivel de representación en la planimetría en escala 1:50. `n- No presentaste las secuencias de uso, las cuales debían identificarse mientras estábamos revisando la entrega.`n- La nota global de tu entrega es 3,7, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 2,7 (rebaja de 1 unidad)."
$ws.Range("G117").Value = "> Comentarios personalizados:`n- Es urgente que revises la manera como estás considerando el entorno que rodea a tu lote y la respuesta que das con relación a este. Esto afecta las posbilidades de iluminar y ventilar de manera adecuada tus espacios.`n- Importante ajustar las proporciones que tienen los diferentes espacios con relación a su uso y optimizar la manera de circular entre estos.`n- Explora alternativas para resolver mejor las cubiertas.`n- La nota global de tu entrega es 2,7."
$ws.Range("G142").Value = "> Comentarios personalizados:`n- Te hizo falta presentar: sección de vía, perspectiva interna y reservar espacio para la sección por fachada. Aunque aparece algo que puede interpretarse como esquemas, debes usar esto de una manera más contundente.`n- Es necesario que eleves el nivel de representación técnica planimétrica y la diagramación de la memoria.`n- No presentaste las secuencias de uso, las cuales debían identificarse mientras estábamos revisando la entrega.`n- La nota global de tu entrega es 3,2."
$ws.Range("G167").Value = "> Comentarios personalizados:`n- Los diferentes espacios que integran tu vivienda están bien ventilados e iluminados. Muy interesante que estés explorando intereses particulares (baños secos y huertas), pues genera un diferencial en tu perfil como arquitecta en coherencia con retos y posibilidades actuales.`n- La manera como se resuelve la cubierta y la proporción de los espacios con relación a su uso, están bien de manera general (es posible hacer ajustes menores en este aspecto).`n- No presentaste perspectiva interna, ni cuadro de áreas.`n- La nota global de tu entrega es 4,2, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 3,2 (rebaja de 1 unidad)."
$ws.Range("G242").Value = "> Comentarios personalizados:`n- El proyecto se está resolviendo de una manera adecuada. Un par de aspectos podrían refinarse mejor: resolución de cubiertas y verificación de topografía (en caso de tener diferencias de nivel, esto debe verse claramente en la representación del proyecto).`n- Explora alternativas para elevar la representación planimétrica de tu proyecto.`n- La nota global de tu entrega es 4,6."
$ws.Range("G267").Value = "> Comentarios personalizados:`n- Un proyecto muy bien desarrollado conforme al nivel de complejidad actual del taller. Estás lista para seguir desarrollándolo a detalle, con una serie de laboratorios que tendremos para el cierre del taller.`n- Solamente vemos que puedes mejorar ligeramente en cuanto a algunos aspectos de representación.`n- La nota global de tu entrega es 4,9."
$ws.Range("G292").Value = "> Comentarios personalizados:`n- Un proyecto que se relaciona de manera adecuada con el entorno que lo rodea.`n- En su interior vemos espacios con buenas posibilidades de iluminarse y ventilarse. En cuanto a áreas de circulación y proporción de espacios con relación a su uso, vas bien y solo vemos algunas posibilidades de refinarse un poco mejor. Para las cubiertas si es necesario explorar más posibilidades pues es el aspecto donde vemos más posibilidades para mejorar.`n- Explora alternativas para elevar el nivel de representación planimétrica.`n- La nota global de tu entrega es 4,6, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 3,6 (rebaja de 1 unidad)."
$ws.Range("G317").Value = "> Comentarios personalizados:`n- Vemos un proyecto que ofrece muy buenas posibilidades de iluminación y ventilación para sus espacios.`n- Es posible explorar con más detalle la manera como se están resolviendo las cubiertas. `n- Te invitamos a explorar alternativas que te permitan relacionarte mejor con el entorno que te rodea, haciéndote preguntas no solo por lo que sucede en el interior de la vivienda, sino en lo que puede aportar tu proyecto para el espacio público de Honda.`n- La manera como se resuelve la circulación al interior de la vivienda, así como la proporción de los espacios con relación a su uso, es donde se pueden hacer los mayores ajustes a la propuesta actual.`n- Es necesario que eleves el nivel de representación, tanto a nivel planimétrico como de diagramación en la memoria + No presentaste perspectiva interna y aunque hay algo relacionado con esquemas, es necesario aprovechar mejor este recurso para presentar tu proyecto.`n- La nota global de tu entrega es 4."
$ws.Range("G342").Value = "> Comentarios personalizados:`n- Encontramos un proyecto con un buen nivel de avance y con relaciones muy claras para el entorno en que se implanta.`n- Es posible refinar la manera como se circula a través de los diferentes espacios. `n- Las cubiertas son en este momento el componente donde se requiere mayor desarrollo y exploración.`n- Explora alternativas para elevar el nivel de representación. `n- La nota global de tu entrega es 4,3."
$ws.Range("G367").Value = "> Comentarios personalizados:`n- Un proyecto muy bien resuelto en sus diferentes aspectos para el nivel de complejidad actual del taller.`n- La representación es donde vemos más posibilidades para mejorar.`n- La nota global de tu entrega es 4,8."
$ws.Range("G392").Value = "> Comentarios personalizados:`n- El proyecto cuenta con un nivel de avance muy bueno para el momento actual del taller. `n- Vemos algunas posibilidades de optimizar un poco mejor la manera como se circula en su interior y la proporción de los espacios con relación al uso que tienen.`n- Te retamos a mejorar tu nivel de representación.`n- La nota global de tu entrega es 4,7."
$ws.Range("G417").Value = "> Comentarios personalizados:`n- De acuerdo a la localización presentada, el proyecto no responde de manera adecuada a las condiciones del entorno en que se implanta (tiene vecinos en 3 costados... no está en una esquina). Esto compromete las posibilidades de iluminar y ventilar de manera adecuada los espacios de alta permanencia. `n- Las cubiertas pueden resolverse de una manera más adecuada.`n- Muy buen nivel de representación técnica y de diagramación.`n- La nota global de tu entrega es 3,9."
$ws.Range("G442").Value = "> Comentarios personalizados:`n- El proyecto se relaciona de manera adecuada con el entorno en que se implanta. También se observa una buena iluminación y ventilación, así como una proporción adecuada según el uso de cada espacio.`n- Es posible refinar la manera como se circula a través de los diferentes espacios.`n- El aspecto que requiere más exploración y resolución es el de la cubierta.`n- Es necesario que eleves el nivel de representación planimétrica y diagramación de la memoria.`n- La nota global de tu entrega es 4,4."
$ws.Range("G467").Value = "> Comentarios personalizados:`n- Es necesario hacer ajustes en el proyecto en cuanto a su relación con el entorno en que se implanta, la manera como se circula a través de los diferentes espacios, la proporción de estos con relación a su uso, ventilación e iluminación, así como la resolución de las cubiertas.`n- Te hizo falta presentar la planta del primer piso, planta de cubiertas, fachada y sección por fachada. Adicionalmente, no presentaste la actividad que estaban desarrollando en paralelo mientras revisábamos la entrega (secuencias de uso)`n- Es importante que explores cómo puedes usar el tiempo de manera más efectiva y de relacionarte con otras personas: esto puede hacer la diferencia en tu proceso de aprendizaje.`n- La nota global de tu entrega es 2,6, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 1,6 (rebaja de 1 unidad)."
$ws.Range("G492").Value = "> Comentarios personalizados:`n- El proyecto está bien resuelto en su generalidad. Identificamos posibilidades de hacer ajustes en la proporción de los espacios con relación a su uso y la manera como se están planteando las cubiertas.`n- Explora alternativas para mejorar la representación planimétrica y diagramar la memoria.`n- La nota global de tu entrega es 4,7."
$ws.Range("G517").Value = "> Comentarios personalizados:`n- El proyecto tiene un nivel de avance adecuado para el momento actual del taller. `n- El aspecto que requiere de un mejor desarrollo está relacionado con la manera como se resuelve la cubierta. La circulación entre los diferentes espacios podría ajustarse un poco.`n- No incorporaste esquemas en la memoria.`n- Es necesario que de cara a la entrega final, eleves el nivel de representación en planimetría y explorar posibilidades de diagramación.`n- La nota global de tu entrega es 4,2."
$ws.Range("G542").Value = "> Comentarios personalizados:`n- Es necesario tener mayor claridad en la manera como te relacionas con el entorno que te rodea y lo que haces con el lote en que implantas el proyecto.`n- El detalle al que tienes desarrollado el proyecto en su interior en este momento, impide ver cómo estás iluminando y ventilando los diferentes espacios o cómo estás resolviendo la cubierta. La cantidad de espacio que destinas para circular y la proporción de los espacios de estancia, requiere de ajustes importantes.`n- No presentaste la planta de cubiertas, ni la sección de vía y esquemas, ni reservaste espacio para la sección por fachada.  encontramos varios de los productos. Tampoco vimos el set de planimetría que requerías presentar para la entrega.`n- Explora maneras de diagramar la memoria en todo su conjunto.`n- Aún estás a tiempo de llegar a un buen puerto en este taller, pero para eso es necesario que explores alternativas de trabajo que te permitan aprovechar mejor el tiempo en clase y ser más efectiva en el trabajo autónomo entre sesiones.`n- La nota global de tu entrega es 2,1."
$ws.Range("G567").Value = "> Comentarios personalizados:`n- Puedes dejar de lado preguntas asociadas a la manera en que tu proyecto se relaciona con el entorno cercano: ya no es necesario hacer grandes ajustes en este aspecto. Ahora es importante de manera prioritaria, resolver la circulación entre los diferentes espacios, de manera eficiente y clara; lo mismo aplica para las cubiertas.`n- Te hace falta reservar el espacio que requiere la sección por fachada.`n- Aprovecha el cierre del taller para elevar tu nivel de representación técnica planimétrica y de diagramación de la memoria.`n- La nota global de tu entrega es 3,9, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 2,9 (rebaja de 1 unidad)."
$ws.Range("G592").Value = "> Comentarios personalizados:`n- Tienes un proyecto desarrollado a un buen nivel para el momento en que nos encontramos. En lo que viene, es posible hacer ajustes menores en la manera como se circula al interior de la vivienda y la proporción que tienen los espacios con relación a su uso.`n- En cuanto a representación logras un buen resultado, con posibilidades de refinar la expresión para una escala 1:50 y la diagramación de la memoria.`n- La nota global de tu entrega es 4,7."
$ws.Range("G617").Value = "> Comentarios personalizados:`n- Cuentas con una localización privilegiada que te permite ofrecer espacialidades más ricas y resolver la envolvente con mayor complejidad. No obstante, parece como si estuvieras trabajando en un lote que tiene muchas restricciones en cuanto a estos aspectos.`n- Al interior de la vivienda es necesario hacer ajustes importantes en la proporción que tienen los espacios con relación a su uso y la eficiencia de las áreas de circulación. Las cubiertas, aunque ya tienen pendientes, no evidencian otras preguntas que ya deberías estar resolviendo: no se trata de una línea inclinada solamente, sino que tienen espesor, debido a una serie de elementos presentes en las mismas.`n- No entregaste el laboratorio formulado para desarrollar mientras evaluábamos la entrega (secuencias de uso). Tampoco presentaste la sección principal, ni la sección de vía, ni reservaste el espacio para la sección por fachada, ni presentaste la perspectiva interna y los esquemas.`n- Hace rato no te asesoras con ninguno de los profes y esto se evidencia en el nivel de desarrollo que tienes en este momento. Aún puedes sacar adelante este proyecto, pero requerirá de tu parte una actitud diferente y un nivel de trabajo mayor al que has invertido a lo largo del semestre.`n- La nota global de tu entrega es 1,7, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 0,7 (rebaja de 1 unidad)."
$ws.Range("G642").Value = "> Comentarios personalizados:`n- El nivel de desarrollo del proyecto es muy bueno y vemos un cuidado importante por la calidad de aquello que presentas.`n- Vemos posibilidades menores de para mejorar la manera como se circular al interior de la vivienda, pero sobre todo, resolver las cubiertas de una manera más adecuada.`n- No encontramos la planimetría en escala 1:250.`n- La nota global de tu entrega es 4,5."
$ws.Range("G667").Value = "> Comentarios personalizados:`n- La manera en que el proyecto se relaciona con el entorno cercano es adecuada y hay una buena respuesta ante las condiciones particulares del lote. Así mismo, la circulación, proporción de espacios según el uso y la manera como se ventilan e iluminan estos, están bien resueltas.`n- Es necesario que mejores tu nivel de representación, tanto en la planimetría técnica, como en la diagramación de la memoria.`n- Te hizo falta presentar la perspectiva interna y los esquemas.`n- La nota global de tu entrega es 4,4."
$ws.Range("G692").Value = "> Comentarios personalizados:`n- El proyecto se relaciona de manera adecuada con el entorno cercano y resuelve de manera coherente con las particularidades del lote. `n- La manera como se circula al interior de la vivienda puede ajustarse mejor, así como la proporción que tienen los diferentes espacios con relación a su uso.`n- Rétate en el nivel que puedes alcanzar en cuanto a representación planimétrica.`n- La nota global de tu entrega es 4,5."
$ws.Range("G717").Value = "> Comentarios personalizados:`n- Proyecto bien resuelto en lo que tiene que ver con el entorno y lo que exige el lote en el que se implanta.`n- La manera como se circula en el interior de la vivienda y el desarrollo de las cubiertas, puede resolverse de una manera más adecuada.`n- Rétata para elevar el nivel que logras en representación técnica planimétrica y de diagramación de la memoria.`n- La nota global de tu entrega es 4,6."
$ws.Range("G742").Value = "> Comentarios personalizados:`n- Es necesario resolver de una manera más adecuada, la relación del proyecto con el entorno cercano y las particularidades del lote.`n- La proporción de algunos espacios con relación a su uso, puede ajustarse un poco mejor.`n- La nota global de tu entrega es 4,8."
$ws.Range("G767").Value = "> Comentarios personalizados:`n- Proyecto con un buen nivel de avance para el momento del taller, tanto por la manera como se relaciona con el entorno, como por su resolución al interior en cuanto a circulacióm proporción de espacios con relación al uso, iluminación y ventilación.`n- El aspecto que requiere de un mejor desarrollo en cuanto a diseño, tiene que ver con la cubierta.`n- Explora alternativas para diagramar mejor la memoria y representar la planimetría en escala 1:50.`n- La nota global de tu entrega es 4,8."
$ws.Range("G792").Value = "> Comentarios personalizados:`n- Es posible hacer ajustes menores en la manera como se relaciona el proyecto con relación al entorno cercano y las particularidades del lote. `n- El aspecto que requiere de mayor atención en este momento en cuanto al diseño, está relacionado con la manera de iluminar y ventilar los espacios de alta permanencia.`n- Así mismo, es posible ajustar mejor la proporción que tienen algunos espacios con relación a su uso, la manera para circular en el interior y el planteamiento de las cubiertas.`n- Es necesario que mejorar tu nivel de representación técnica planimétrica, pero sobre todo, la manera como diagramas diferentes productos dentro de una memoria.`n- Te hace falta presentar la sección principal, la sección de vía y reservar el espacio para la sección por fachada. Tampoco hay esquemas ni cuadro de áreas. No entregaste el laboratorio formulado para hacerse mientras evaluábamos la entrega (secuencias de uso).`n- Hace varias sesiones que no conversas con ninguno de los profes y esto afecta el nivel de avance que tiene tu proyecto.`n- La nota global de tu entrega es 3, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 2 (rebaja de 1 unidad)."
$ws.Range("G817").Value = "> Comentarios personalizados:`n- El proyecto puede relacionarse de una manera más adecuada con el entorno que lo rodea y las particularidades del lote en que se encuentra.`n- De manera prioritaria, resuelve mejor la proporción de algunos espacios con relación a su uso, la iluminación y ventilación de estos, así como lo que estás planteando para la cubierta.`n- Te hace falta reservar el espacio para la sección por fachada.`n- Eleva el nivel de representación técnica planimétrica y la manera de diagramar diferentes productos en una memoria.`n- La nota global de tu entrega es 4,2, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 3,2 (rebaja de 1 unidad)."
$ws.Range("G842").Value = "> Comentarios personalizados:`n- El proyecto está bien resuelto para el momento actual del taller, tanto en su relación con el entorno y la interpretación del lote, así como en su circulación, proporción, iluminación y ventilación de espacios.`n- No presentaste el set de planos y la diagramación de la memoria puede mejorar + utiliza esquemas que te permitan comunicar ideas clave para entender tu proyecto.`n- La nota global de tu entrega es 3,9."
$ws.Range("G867").Value = "> Comentarios personalizados:`n- Es necesario de manera prioritaria, ajustar la manera de circular al interior de la vivienda.`n- De otra parte, la proporción de algunos espacios puede ajustarse mejor de acuerdo a su uso. Así mismo, es necesario explorar mejor la manera como se resuelven las cubiertas.`n- Tienes un margen de mejora importante en la representación técnica planimétrica en escala 1:50 + ensaya más alternativas para diagramar el contenido en una memoria.`n- La nota global de tu entrega es 4,2."
$ws.Range("G892").Value = "> Comentarios personalizados:`n- El proyecto se relaciona de manera adecuada con relación al entorno cercano y responde bien las particularidades del lote. Así mismo, el planteamiento de cubiertas actual está muy bien para el momento actual del taller.`n- Es posible resolver un poco mejor la manera como se circula en su interior y la proporción que tienen algunos espacios con relación a su uso.`n- El área en el que tienes mayor margen de mejora, tiene que ver con el nivel de representación, tanto a nivel técnico planimétrico, como de diagramación de la memoria.`n- La nota global de tu entrega es 4,2."
$ws.Range("G917").Value = "> Comentarios personalizados:`n- Es necesario resolver de una manera más adecuada, la relación del proyecto con el entorno cercano y las particularidades del lote.`n- Es importante explorar alternativas para iluminar y ventilar mejor algunos espacios. Así mismo, la proporción de algunos espacios con relación a su uso, puede ajustarse un poco mejor. Lo mismo aplica para la manera como circulas a través de estos.`n- La nota global de tu entrega es 4,4."
$ws.Range("G942").Value = "> Comentarios personalizados:`n- El proyecto se relaciona de manera adecuada con su entorno y responde a las particularidades del lote en que se implanta + La diagramación de la memoria es muy buena.`n- En cuanto al diseño arquitectónico, es posible ajustar las proporciones de algunos espacios con relación a su uso, así como la manera como se circula a través de estos. En unos cuantos, es posible mejorar su ventilación e iluminación.`n- Aunque presentaste a tiempo todos los productos requeridos para la entrega a través de la memoria, el set de planos nunca apareció.`n- La nota global de tu entrega es 4."
$ws.Range("G967").Value = "> Comentarios personalizados:`n- De manera prioritaria, debes explorar alternativas para resolver las cubiertas de los diferentes volúmenes. Una manera que puede ayudarte mucho para esto, es aplicar una lógica geométrica similar a la que usaste para la planta del primer nivel (hexágonos).`n- Es posible hacer algunos ajustes en la manera como se circula en el interior de la vivienda y la proporción que tienen algunos espacios con relación a su uso.`n- Es necesario que mejores mucho en la diagramación de la memoria y el nivel de representación planimétrica.`n- No reservaste espacio para la sección por fachada, ni presentaste esquemas ni cuadro de áreas.`n- La nota global de tu entrega es 3,8, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 2,8 (rebaja de 1 unidad)."
$ws.Range("G992").Value = "> Comentarios personalizados:`n- El proyecto se relaciona de manera adecuada con relación a su entorno cercano y las posibilidades que ofrece el lote en que se implanta. La ventilación e iluminación de los espacios está bien resuelta.`n- El aspecto que puede desarrollarse mejor en cuanto al diseño arquitectónico, tiene que ver con la manera de circular al interior de la vivienda y el planteamiento de las cubiertas. La proporción de algunos espacios puede ajustarse mejor con relación a su uso.`n- Explora alternativas para mejorar la representación técnica planimétrica.`n- La nota global de tu entrega es 4,3."
$ws.Range("G1017").Value = "> Comentarios personalizados:`n- Conviene explorar alternativas en cuanto a la manera en que el proyecto responde a las condiciones del lote y el entorno cercano.`n- Al interior de la vivienda es importante revisar la proporción de algunos espacios con relación a su uso y la ventilación e iluminación de algunos de estos.`n- Te hace falta incorporar esquemas en la memoria, que te permitan comunicar de manera efectiva ideas potentes de tu diseño.`n- La nota global de tu entrega es 4,6."
$ws.Range("G1042").Value = "> Comentarios personalizados:`n- De manera prioritaria, explora alternativas para resolver mejor la cubierta del proyecto.`n- Hay algunos espacios de alta permanencia que pueden iluminarse y ventilarse mejor. Así mismo, es posible ajustar la proporción de algunos de éstos, en función del uso que tendrán.`n- El nivel de representación planimétrica tiene mucho margen de mejora, así como la diagramación de la memoria. Nos sorprendió lo que presentaste para esta entrega, pues has demostrado muchas capacidades a lo largo del semestre.`n- No presentaste el laboratorio planteado para desarrollarse en paralelo a la evaluación de la entrega (secuencias de uso).`n- La nota global de tu entrega es 3,4, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 2,4 (rebaja de 1 unidad)."
$ws.Range("G1067").Value = "> Comentarios personalizados:`n- La manera como se relaciona el proyecto con el entorno puede mejorar. Así mismo, hay ajustes que se pueden hacer en cuanto a circulación, iluminación y ventilación, cubiertas y la proporción de los espacios con relación a su uso.`n- No presentaste la perspectiva interna, ni los esquemas. Tampoco reservaste el espacio requerido para la sección por fachada. Adicionalmente, no presentaste el laboratorio planteado para desarrollarse mientras evaluábamos la entrega (secuencias de uso).`n- Es muy importante que explores alternativas para diagramar diferentes productos dentro de una memoria. La representación técnica planimétrica también debe subir de nivel.`n- La nota global de tu entrega es 3,4, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 2,4 (rebaja de 1 unidad)."
$ws.Range("G1092").Value = "> Comentarios personalizados:`n- El proyecto se relaciona de manera adecuada con el entorno y responde coherentemente a las particularidades del lote en que se implanta.`n- En su interior, hay un buen desarrollo, con algunas opciones para refinar la proporción de algunos espacios con relación al uso y el planteamiento de la cubierta.`n- La representación técnica planimétrica debe subir de nivel para lo que se espera en este nivel de formación. Así mismo, la diagramación de la memoria, tiene margen de mejora importante.`n- No presentaste fachada, ni esquemas.`n- La nota global de tu entrega es 4,1, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 3,1 (rebaja de 1 unidad)."
$ws.Range("G1117").Value = "> Comentarios personalizados:`n- El proyecto se relaciona de manera adecuada con su entorno cercano y las particularidades del lote.`n- En donde más puedes mejorar, tiene que ver con la representación, sobre todo en cuanto a diagramación de la memoria.`n- En cuanto al diseño arquitectónico, es prioritario resolver de una manera más adecuada las cubiertas.`n- La circulación y proporción de espacios según su uso, aunque van por buen camino, pueden refinarse un poco mejor.`n- Hace falta presentar fachada y esquemas`n- La nota global de tu entrega es 4,1."
$ws.Range("G1142").Value = "> Comentarios personalizados:`n- El proyecto está bien resuelto en su interior, para el momento actual del taller (circulación, proporción de espacios, iluminación y ventilación de estos). El único aspecto que se puede mejorar, está relacionado con la cubierta.`n- Donde requieres ocuparte más en este momento, es en la representación técnica planimétrica y la diagramación de la memoria.`n- La nota global de tu entrega es 4,4, pero al haber entregado productos de manera extemporánea al límite anunciado, se computa en 3,4 (rebaja de 1 unidad)."

# --- Row 192-196 (Fabian Enrique Arevalo Romero): new comment text + grade updates ---
$ws.Range("G192").Value = "> Comentarios personalizados:`n- Es necesario presentar planta de localización y no una aerofotografía en perspectiva, fachada y no 3D en el que se ven las fachadas. Necesitas de esos dibujos en vistas ortogonales, de tal manera que se puedan evaluar y ajustar detalles del proyecto con relación al entorno cercano, entre otros procesos.`n- Al estar en una manzana tan larga, adquieres tácitamente una responsabilidad con relación al espacio público de Honda: ¿qué le vas a ofrecer a  los peatones que circulen por ese andén?`n- Pese a haber presentado la entrega en una sesión posterior, te hicieron falta varios de los productos requeridos. Para el cierre del taller es fundamental que revises la manera como estás trabajando y así poder alcanzar los resultados de aprendizaje esperados para el curso.`n- La nota global de tu entrega es 1,5."
$ws.Range("C193").Value = "'1.6"
$ws.Range("C194").Value = "'0.7"
$ws.Range("C195").Value = "'1.7"
$ws.Range("C196").Value = "'2.2"

# --- Row 217-221 (George Five Quinto Moya): new comment text + rubric corrections ---
$ws.Range("G217").Value = "> Comentarios personalizados:`n- Considera los comentarios realizados en la sesión posterior al momento de la entrega. Con los productos presentados en ese momento, se aclaran algunas dudas que teníamos y se abren otras que ya conoces. `n- La nota global de tu entrega es 2,9."
$ws.Range("C218").Value = "'1.9"
$ws.Range("C219").Value = "'2.2"
$ws.Range("D219").Value = "C"
$ws.Range("F219").Value = "E4.Rep_C"
$ws.Range("G219").Value = "La representación y comunicación del proyecto no es clara. El lenguaje arquitectónico y la estética empleada, no corresponde al nivel de formación."
$ws.Range("C220").Value = "'3.3"
$ws.Range("C221").Value = "'4.3"
$ws.Range("D221").Value = "A"
$ws.Range("F221").Value = "E4.Ser_A"
$ws.Range("G221").Value = "Desarrolla el proyecto con una actitud crítica, respetuosa, empática, propositiva y cooperativa. Demuestra habilidades de comunicación efectiva, pensamiento crítico y creatividad de manera sobresaliente."
